$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cited_by_count for row 2 (M2): 7 -> 8
# Force the cell to remain text (matches the original inlineStr "7"),
# instead of Excel auto-converting the numeric-looking string to a number.
$m2 = $ws.Range("M2")
$m2.NumberFormat = "@"
$m2.Value2 = "8"
$m2.Style = "Normal"

# Swap the contents of row 9 and row 10 for columns A, B, C, D, H, P
$cols = @("A", "B", "C", "D", "H", "P")

foreach ($col in $cols) {
    $cell9 = $ws.Range($col + "9")
    $cell10 = $ws.Range($col + "10")
    $temp = $cell9.Value2
    $cell9.Value2 = $cell10.Value2
    $cell10.Value2 = $temp
}
